$d = $word.ActiveDocument

# --- 1) "문제 정의 (Problem Definition)" heading: bump bold heading font size 12pt -> 14pt ---
# (affects both the paragraph-mark run props and the two existing runs)
$pProblem = $d.Paragraphs.Item(1)
$pProblem.Range.Font.Size = 14
$pProblem.Range.Font.SizeBi = 14

# --- 2) "아이디어 생성" heading: bump bold heading font size 12pt -> 14pt ---
$pIdea = $d.Paragraphs.Item(4)
$pIdea.Range.Font.Size = 14
$pIdea.Range.Font.SizeBi = 14

# --- 3) Append "(Idea Generation)" after "아이디어 생성", split across two runs so the
#        "(" keeps the eastAsia hint (it follows Korean text) while "Idea Generation)"
#        (pure Latin) does not - matching how the existing "문제 정의 (Problem Definition)"
#        heading above is already split. We rebuild the whole paragraph (mark included) via
#        InsertXML so we can specify the run boundaries / rFonts precisely; the paragraph's
#        own identity (paraId/rsids) and pPr are supplied to keep them stable, and its pPr
#        is only actually honored for the size we already set on the mark above.
$pIdea2 = $d.Paragraphs.Item(4)
$fullRange = $d.Range($pIdea2.Range.Start, $pIdea2.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="144420B6" w14:textId="77777777" w:rsidR="00A41030" w:rsidRPr="00F97BB1" w:rsidRDefault="00A41030" w:rsidP="00A41030"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00F97BB1"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>아이디어 생성</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Idea Generation)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$fullRange.InsertXML($xml)

Write-Output "Applied InClass idea-generation heading update."
